# Cotações atualizadas - 2025-10-10
# Append a new row (36) with the latest quotes, following the same
# pattern as the existing rows (date serial in column A, comma-decimal
# text values in columns B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 36

# Column A: date serial value, formatted the same way as the previous row.
$ws.Range("A$newRow").Value() = 45940
$ws.Range("A$newRow").NumberFormat = $ws.Range("A35").NumberFormat

# Columns B:E: quote values stored as text (comma decimal separator).
$ws.Range("B$newRow").Value() = "21,6987"
$ws.Range("C$newRow").Value() = "15,6648"
$ws.Range("D$newRow").Value() = "15,4517"
$ws.Range("E$newRow").Value() = "15,4517"
